$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New lower_snake_case headers
$ws.Range("A1").Value = "nome_completo"
$ws.Range("B1").Value = "nome_curto"
$ws.Range("C1").Value = "cnpj"

# Column B (nome_curto) now mirrors column A (nome_completo) for every data row,
# and column A gets underlined to flag it.
for ($r = 2; $r -le 9; $r++) {
    $nome = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $nome
    $ws.Cells.Item($r, 1).Font.Underline = $true
}

# Mirror the selection left behind in the saved file (A2:A9, active cell A9)
$ws.Range("A2:A9").Select()
$excel.ActiveCell = $ws.Range("A9")
